$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.018.59"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.643.04"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.14"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "1.755.91"
$ws.Range("E12").Value = "  +7.35%  "
$ws.Range("D13").Value = "1.871.10"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.28"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.41"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "26.110.74"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.52"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.22"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.131"
$ws.Range("E24").Value = "  +4.71%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "143.31"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0497"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("E34").Value = "  -2.86%  "
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").Value = "1.133.13"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.40"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "1.781.02"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "0.0₆0118"
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.74"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  +2.52%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0959"
$ws.Range("E51").Value = "  -0.08%  "
